# 2024.11.30 - Save Final Output File
# Updates the "Correlation Sample" sheet (B2:P16) with recomputed correlation values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Correlation Sample")

$values = @(
  @(1,0.59,0.43,0.4,0.46,0.29,0.15,0.36,0.57,0.49,0.27,0.65,0.33,0.47,0.16),
  @(0.59,1,0.43,0.27,0.55,0.15,0.17,0.43,0.5600000000000001,0.57,0.05,0.58,0.12,0.3,0.05),
  @(0.43,0.43,1,0.36,0.44,0.41,0.28,0.4,0.49,0.52,0.17,0.33,0.21,0.35,0.15),
  @(0.4,0.27,0.36,1,0.33,0.55,0.17,0.27,0.32,0.34,0.33,0.25,0.35,0.54,0.42),
  @(0.46,0.55,0.44,0.33,1,0.35,0.13,0.41,0.55,0.51,0.11,0.42,0.21,0.36,0.22),
  @(0.29,0.15,0.41,0.55,0.35,1,0.09,0.24,0.24,0.29,0.17,0.19,0.29,0.48,0.52),
  @(0.15,0.17,0.28,0.17,0.13,0.09,1,-0,0.33,0.16,0.11,0.03,0.1,0.16,0.14),
  @(0.36,0.43,0.4,0.27,0.41,0.24,-0,1,0.47,0.5,0.25,0.3,0.06,0.28,-0.02),
  @(0.57,0.5600000000000001,0.49,0.32,0.55,0.24,0.33,0.47,1,0.5600000000000001,0.25,0.45,0.26,0.34,0.1),
  @(0.49,0.57,0.52,0.34,0.51,0.29,0.16,0.5,0.5600000000000001,1,0.03,0.4,0.03,0.27,0.01),
  @(0.27,0.05,0.17,0.33,0.11,0.17,0.11,0.25,0.25,0.03,1,-0.01,0.27,0.31,0.11),
  @(0.65,0.58,0.33,0.25,0.42,0.19,0.03,0.3,0.45,0.4,-0.01,1,0.21,0.29,0.06),
  @(0.33,0.12,0.21,0.35,0.21,0.29,0.1,0.06,0.26,0.03,0.27,0.21,1,0.28,0.22),
  @(0.47,0.3,0.35,0.54,0.36,0.48,0.16,0.28,0.34,0.27,0.31,0.29,0.28,1,0.44),
  @(0.16,0.05,0.15,0.42,0.22,0.52,0.14,-0.02,0.1,0.01,0.11,0.06,0.22,0.44,1)
)

$rowCount = $values.Count
$colCount = $values[0].Count
$arr = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
  for ($j = 0; $j -lt $colCount; $j++) {
    $arr[$i,$j] = $values[$i][$j]
  }
}

$ws.Range("B2:P16").Value = $arr

Write-Output "Updated Correlation Sample B2:P16 ($rowCount x $colCount)"
